$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1650
$ws1.Range("F3").Value = 9061
$ws1.Range("F6").Value = 696
$ws1.Range("F7").Value = 956
$ws1.Range("F8").Value = 191
$ws1.Range("F11").Value = 5675
$ws1.Range("F15").Value = 4363
$ws1.Range("F19").Value = 18
$ws1.Range("F24").Value = 2704

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1650
$ws4.Range("F3").Value = 9061
$ws4.Range("F5").Value = 14
$ws4.Range("F7").Value = 696
$ws4.Range("F8").Value = 957
$ws4.Range("F9").Value = 191
$ws4.Range("F12").Value = 5675
$ws4.Range("F16").Value = 4363
$ws4.Range("F20").Value = 18
$ws4.Range("F25").Value = 2704
